$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BackLog")

# --- Update PROGRESSO column: "Andamento" -> "Concluído" for the first four tasks ---
$ws.Range("F3").Value = "Concluído"
$ws.Range("F4").Value = "Concluído"
$ws.Range("F5").Value = "Concluído"
$ws.Range("F6").Value = "Concluído"

# --- Add a new backlog entry (E11) in row 13, matching the formatting of row 12 ---
$ws.Range("B12:F12").Copy($ws.Range("B13:F13"))
$excel.CutCopyMode = $false

$ws.Range("B13").Value = "E11"
$ws.Range("C13").Value = "Fazer BPMN da aplicação"
$ws.Range("D13").Value = "Importante"
$ws.Range("E13").Value = 5
$ws.Range("F13").Value = "Concluído"

# --- Move the active cell selection ---
$ws.Range("H18").Select()
